# Apply "Added new GA 3.2 prototype" / "Added Course Sections" edit to schedule.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet
$ws.Name = "Class Schedule"

# 2. Strip the old bold+bordered header style from the weekday header row
#    and from the time-of-day column so they fall back to the default style.
$ws.Range("B1:G1").ClearFormats()
$ws.Range("A2:A29").ClearFormats()

# 3. Add the new "Time" column header
$ws.Range("A1").Value = "Time"

# 4. Wipe the old (empty placeholder + stale course) cells in the grid body
$ws.Range("B2:G29").ClearContents()

# 5. Populate the new course sections, centre them, and merge their time spans
$sections = @(
    @{ Range = "B4:B6";   Text = "CS13 - CS101 (Room 1)" },
    @{ Range = "E4:E6";   Text = "CS13 - CS101 (Room 1)" },
    @{ Range = "C6:C11";  Text = "CS13 - CS102 (Room 1)" },
    @{ Range = "F6:F11";  Text = "CS13 - CS102 (Room 1)" },
    @{ Range = "B8:B10";  Text = "CS12 - CS101 (Room 1)" },
    @{ Range = "E8:E10";  Text = "CS12 - CS101 (Room 1)" },
    @{ Range = "D10:D19"; Text = "CS12 - CS103 (Room 2)" },
    @{ Range = "C14:C19"; Text = "CS11 - CS102 (Room 3)" },
    @{ Range = "F14:F19"; Text = "CS11 - CS102 (Room 3)" },
    @{ Range = "D20:D29"; Text = "CS13 - CS103 (Room 4)" }
)

foreach ($section in $sections) {
    $rng = $ws.Range($section.Range)
    $rng.Value = $section.Text
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4108
    $rng.Merge()
}
